$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Per-row Price (D) and Volume(1h) (E) updates.
# D values that look like plain decimals need NumberFormat "@" forced
# before the assignment (else Excel auto-converts them to numbers),
# then the style is reset back to Normal so no stray number format
# sticks to the cell (matches original default-style inline strings).

$ws.Range("D2").Value = '26.231.60'
$ws.Range("E2").Value = '  +0.24%  '

$ws.Range("D3").Value = '1.600.29'
$ws.Range("E3").Value = '  -0.21%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.10'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.02%  '

$ws.Range("E6").Value = '  -0.13%  '

$ws.Range("E7").Value = '  +0.47%  '

$ws.Range("E8").Value = '  -0.02%  '

$ws.Range("E9").Value = '  -0.62%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.08'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.08%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0813'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.26%  '

$ws.Range("D12").Value = '1.822.23'
$ws.Range("E12").Value = '  -0.29%  '

$ws.Range("D13").Value = '1.598.16'
$ws.Range("E13").Value = '  +0.04%  '

$ws.Range("E14").Value = '  +0.86%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.512'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.82%  '

$ws.Range("D16").Value = '26.205.80'
$ws.Range("E16").Value = '  +0.14%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.19'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.99%  '

$ws.Range("D18").Value = '0.0₃0729'
$ws.Range("E18").Value = '  +0.19%  '

$ws.Range("E19").Value = '  -0.08%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '201.95'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.77%  '

$ws.Range("E21").Value = '  +0.91%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.26'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.49%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.99'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.12%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '144.48'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.69%  '

$ws.Range("E26").Value = '  -0.06%  '

$ws.Range("E27").Value = '  -7.49%  '

$ws.Range("E28").Value = '  +0.21%  '

$ws.Range("E29").Value = '  +1.08%  '

$ws.Range("E30").Value = '  +3.74%  '

$ws.Range("E31").Value = '  -0.41%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.16'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.91%  '

$ws.Range("E33").Value = '  -3.12%  '

$ws.Range("E34").Value = '  +2.81%  '

$ws.Range("E35").Value = '  -1.14%  '

$ws.Range("D36").Value = '1.154.33'
$ws.Range("E36").Value = '  +4.23%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0165'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +8.15%  '

$ws.Range("E39").Value = '  -1.00%  '

$ws.Range("E40").Value = '  +0.62%  '

$ws.Range("E41").Value = '  -0.67%  '

$ws.Range("E42").Value = '  +0.42%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.21'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.83%  '

$ws.Range("D44").Value = '1.738.31'
$ws.Range("E44").Value = '  -0.08%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '91.84'
$ws.Range("D45").Style = "Normal"

$ws.Range("E46").Value = '  -2.51%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '54.09'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.05%  '

$ws.Range("E48").Value = '  -0.68%  '

# Rows 49-51 were reshuffled: BabyDogeCoin moves up to row 49,
# pushing Mantle to row 50 and USDD to row 51. Coin/link identity
# travels with each coin; price/volume get refreshed values for
# their new row.
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = '0.0₇0970'
$ws.Range("E49").Value = '  -6.72%  '

$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.407'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.57%  '

$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.00'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.10%  '
